$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Jours de Ramasse S1" header (column E) to "Jours de Ramasse"
$ws.Range("E1").Value = "Jours de Ramasse"

# Remove the "Jours de Ramasse S2" column (F) entirely; this shifts the
# "Poids par ramasse(kg)" column (G) left into the F position, carrying
# over its style/width/values automatically.
$ws.Columns("F").Delete()
